# Update coffee stock inventory quantities (column E) and refresh the
# stock-level highlight (fill) on each affected row so it keeps matching
# the new quantity: 0 -> red (style 4), 1-2 -> yellow (style 5),
# 3+ -> no fill (default style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Donor rows already carrying each of the three stock-level styles so we
# reuse the existing style indices instead of minting new ones.
$RedStyleDonor    = "A2:E2"    # qty = 0  -> red fill
$YellowStyleDonor = "A13:E13"  # qty = 1-2 -> yellow fill
$NoneStyleDonor   = "A11:E11"  # qty >= 3 -> no fill

function Set-RowStyle($rowNum, $donorRange) {
    $ws.Range($donorRange).Copy()
    $ws.Range("A$($rowNum):E$($rowNum)").PasteSpecial(-4122)  # xlPasteFormats
}

# Row 7: Buna Blend 500g Beans -> out of stock (2 -> 0), red fill
Set-RowStyle 7 $RedStyleDonor
$ws.Range("E7").Value = 0

# Row 9: Buna Blend 1kg Beans -> restocked a bit (1 -> 2), stays yellow
$ws.Range("E9").Value = 2

# Row 14: Koke 250g Ground -> well stocked now (1 -> 3), no fill
Set-RowStyle 14 $NoneStyleDonor
$ws.Range("E14").Value = 3

# Row 15: Koke 250g Beans -> more stock (5 -> 7), stays unstyled
$ws.Range("E15").Value = 7

# Row 22: Sidamo 250g Ground -> restocked a bit (1 -> 2), stays yellow
$ws.Range("E22").Value = 2

# Row 23: Sidamo 250g Beans -> back in stock (0 -> 2), yellow fill
Set-RowStyle 23 $YellowStyleDonor
$ws.Range("E23").Value = 2

# Row 30: Yirgacheffe 250g Ground -> out of stock (2 -> 0), red fill
Set-RowStyle 30 $RedStyleDonor
$ws.Range("E30").Value = 0

# Row 58: Glory to Ukraine 250g Ground -> stock dropped (2 -> 1), stays yellow
$ws.Range("E58").Value = 1

# Row 68: As ir Tu 250g Ground -> back in stock (0 -> 2), yellow fill
Set-RowStyle 68 $YellowStyleDonor
$ws.Range("E68").Value = 2

# Row 69: As ir Tu 250g Beans -> back in stock (0 -> 1), yellow fill
Set-RowStyle 69 $YellowStyleDonor
$ws.Range("E69").Value = 1

# Row 74: Upendo Africa 250g Ground -> stock dropped (5 -> 3), stays unstyled
$ws.Range("E74").Value = 3
